# Scheduled market-data refresh: update cached currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ figures
# across all class sheets to the latest pulled values.

$edits = @(
    @{Sheet="ALC"; Cell="H2"; Value=124},
    @{Sheet="ALC"; Cell="I2"; Value=124},
    @{Sheet="ALC"; Cell="J2"; Value=0},
    @{Sheet="ALC"; Cell="K2"; Value=124},
    @{Sheet="ALC"; Cell="L2"; Value=0},
    @{Sheet="ALC"; Cell="M2"; Value=-11},
    @{Sheet="ALC"; Cell="N2"; Value=$null},
    @{Sheet="ALC"; Cell="H62"; Value=5649.9},
    @{Sheet="ALC"; Cell="I62"; Value=1624.75},
    @{Sheet="ALC"; Cell="K62"; Value=1624.75},
    @{Sheet="ALC"; Cell="M62"; Value=-1000.75},
    @{Sheet="ALC"; Cell="H65"; Value=5649.9},
    @{Sheet="ALC"; Cell="I65"; Value=1624.75},
    @{Sheet="ALC"; Cell="K65"; Value=8123.75},
    @{Sheet="ALC"; Cell="M65"; Value=-5003.75},
    @{Sheet="ALC"; Cell="H100"; Value=1228.5454},
    @{Sheet="ALC"; Cell="I100"; Value=783.25},
    @{Sheet="ALC"; Cell="J100"; Value=2416},
    @{Sheet="ALC"; Cell="K100"; Value=783.25},
    @{Sheet="ALC"; Cell="L100"; Value=2416},
    @{Sheet="ALC"; Cell="M100"; Value=-242.25},
    @{Sheet="ALC"; Cell="N100"; Value=-3498},
    @{Sheet="ALC"; Cell="H112"; Value=2162.3157},
    @{Sheet="ALC"; Cell="J112"; Value=2182},
    @{Sheet="ALC"; Cell="L112"; Value=6546},
    @{Sheet="ALC"; Cell="N112"; Value=-8762},
    @{Sheet="ALC"; Cell="H135"; Value=1326},
    @{Sheet="ALC"; Cell="I135"; Value=1162.6522},
    @{Sheet="ALC"; Cell="J135"; Value=1795.625},
    @{Sheet="ALC"; Cell="K135"; Value=10463.8698},
    @{Sheet="ALC"; Cell="L135"; Value=16160.625},
    @{Sheet="ALC"; Cell="M135"; Value=-7928.8698},
    @{Sheet="ALC"; Cell="N135"; Value=-21230.625},
    @{Sheet="ALC"; Cell="H137"; Value=2577.2},
    @{Sheet="ALC"; Cell="I137"; Value=2051.35},
    @{Sheet="ALC"; Cell="K137"; Value=6154.049999999999},
    @{Sheet="ALC"; Cell="M137"; Value=-3604.049999999999},
    @{Sheet="ALC"; Cell="H138"; Value=3962.9048},
    @{Sheet="ALC"; Cell="J138"; Value=3062.889},
    @{Sheet="ALC"; Cell="L138"; Value=9188.667000000001},
    @{Sheet="ALC"; Cell="N138"; Value=-19468.667},
    @{Sheet="ALC"; Cell="H141"; Value=7239},
    @{Sheet="ALC"; Cell="I141"; Value=5298.75},
    @{Sheet="ALC"; Cell="J141"; Value=15000},
    @{Sheet="ALC"; Cell="K141"; Value=15896.25},
    @{Sheet="ALC"; Cell="L141"; Value=45000},
    @{Sheet="ALC"; Cell="M141"; Value=-10716.25},
    @{Sheet="ALC"; Cell="N141"; Value=-55360},
    @{Sheet="ARM"; Cell="H74"; Value=2329.6667},
    @{Sheet="ARM"; Cell="I74"; Value=1199.8125},
    @{Sheet="ARM"; Cell="J74"; Value=5945.2},
    @{Sheet="ARM"; Cell="K74"; Value=1199.8125},
    @{Sheet="ARM"; Cell="L74"; Value=5945.2},
    @{Sheet="ARM"; Cell="M74"; Value=-325.8125},
    @{Sheet="ARM"; Cell="N74"; Value=-7693.2},
    @{Sheet="ARM"; Cell="H77"; Value=2329.6667},
    @{Sheet="ARM"; Cell="I77"; Value=1199.8125},
    @{Sheet="ARM"; Cell="J77"; Value=5945.2},
    @{Sheet="ARM"; Cell="K77"; Value=5999.0625},
    @{Sheet="ARM"; Cell="L77"; Value=29726},
    @{Sheet="ARM"; Cell="M77"; Value=-1631.0625},
    @{Sheet="ARM"; Cell="N77"; Value=-38462},
    @{Sheet="ARM"; Cell="H102"; Value=3000},
    @{Sheet="ARM"; Cell="I102"; Value=0},
    @{Sheet="ARM"; Cell="J102"; Value=3000},
    @{Sheet="ARM"; Cell="K102"; Value=0},
    @{Sheet="ARM"; Cell="L102"; Value=3000},
    @{Sheet="ARM"; Cell="M102"; Value=$null},
    @{Sheet="ARM"; Cell="N102"; Value=-6244},
    @{Sheet="ARM"; Cell="H132"; Value=4003.4546},
    @{Sheet="ARM"; Cell="I132"; Value=3893.2222},
    @{Sheet="ARM"; Cell="K132"; Value=11679.6666},
    @{Sheet="ARM"; Cell="M132"; Value=-9149.6666},
    @{Sheet="BSM"; Cell="H20"; Value=4128.45},
    @{Sheet="BSM"; Cell="I20"; Value=4965.5454},
    @{Sheet="BSM"; Cell="K20"; Value=4965.5454},
    @{Sheet="BSM"; Cell="M20"; Value=-4718.5454},
    @{Sheet="BSM"; Cell="H105"; Value=4561.3},
    @{Sheet="BSM"; Cell="I105"; Value=4734.778},
    @{Sheet="BSM"; Cell="K105"; Value=4734.778},
    @{Sheet="BSM"; Cell="M105"; Value=-2987.778},
    @{Sheet="BSM"; Cell="H122"; Value=0},
    @{Sheet="BSM"; Cell="J122"; Value=0},
    @{Sheet="BSM"; Cell="L122"; Value=0},
    @{Sheet="BSM"; Cell="N122"; Value=$null},
    @{Sheet="BSM"; Cell="H134"; Value=2459.158},
    @{Sheet="BSM"; Cell="I134"; Value=2115.7334},
    @{Sheet="BSM"; Cell="K134"; Value=6347.2002},
    @{Sheet="BSM"; Cell="M134"; Value=-3812.2002},
    @{Sheet="CRP"; Cell="H22"; Value=374.66666},
    @{Sheet="CRP"; Cell="I22"; Value=366},
    @{Sheet="CRP"; Cell="J22"; Value=383.33334},
    @{Sheet="CRP"; Cell="K22"; Value=366},
    @{Sheet="CRP"; Cell="L22"; Value=383.33334},
    @{Sheet="CRP"; Cell="M22"; Value=-16},
    @{Sheet="CRP"; Cell="N22"; Value=-1083.33334},
    @{Sheet="CRP"; Cell="H31"; Value=3618.724},
    @{Sheet="CRP"; Cell="J31"; Value=4818.385},
    @{Sheet="CRP"; Cell="L31"; Value=4818.385},
    @{Sheet="CRP"; Cell="N31"; Value=-5408.385},
    @{Sheet="CRP"; Cell="H34"; Value=3618.724},
    @{Sheet="CRP"; Cell="J34"; Value=4818.385},
    @{Sheet="CRP"; Cell="L34"; Value=4818.385},
    @{Sheet="CRP"; Cell="N34"; Value=-5222.385},
    @{Sheet="CRP"; Cell="H99"; Value=11018.742},
    @{Sheet="CRP"; Cell="J99"; Value=14609.8},
    @{Sheet="CRP"; Cell="L99"; Value=14609.8},
    @{Sheet="CRP"; Cell="N99"; Value=-17605.8},
    @{Sheet="CRP"; Cell="H122"; Value=2499.6843},
    @{Sheet="CRP"; Cell="I122"; Value=2499.6843},
    @{Sheet="CRP"; Cell="K122"; Value=7499.0529},
    @{Sheet="CRP"; Cell="M122"; Value=-5049.0529},
    @{Sheet="CRP"; Cell="H126"; Value=11018.742},
    @{Sheet="CRP"; Cell="J126"; Value=14609.8},
    @{Sheet="CRP"; Cell="L126"; Value=43829.39999999999},
    @{Sheet="CRP"; Cell="N126"; Value=-48769.39999999999},
    @{Sheet="CRP"; Cell="H132"; Value=2470.9565},
    @{Sheet="CRP"; Cell="I132"; Value=2501.0908},
    @{Sheet="CRP"; Cell="J132"; Value=1808},
    @{Sheet="CRP"; Cell="K132"; Value=7503.2724},
    @{Sheet="CRP"; Cell="L132"; Value=5424},
    @{Sheet="CRP"; Cell="M132"; Value=-4973.2724},
    @{Sheet="CRP"; Cell="N132"; Value=-10484},
    @{Sheet="CRP"; Cell="H134"; Value=2037.8182},
    @{Sheet="CRP"; Cell="I134"; Value=1007.9231},
    @{Sheet="CRP"; Cell="K134"; Value=3023.7693},
    @{Sheet="CRP"; Cell="M134"; Value=-488.7692999999999},
    @{Sheet="CUL"; Cell="H12"; Value=541.96155},
    @{Sheet="CUL"; Cell="I12"; Value=508},
    @{Sheet="CUL"; Cell="J12"; Value=554.4737},
    @{Sheet="CUL"; Cell="K12"; Value=1524},
    @{Sheet="CUL"; Cell="L12"; Value=1663.4211},
    @{Sheet="CUL"; Cell="M12"; Value=-1351},
    @{Sheet="CUL"; Cell="N12"; Value=-2009.4211},
    @{Sheet="CUL"; Cell="H117"; Value=2999.2856},
    @{Sheet="CUL"; Cell="J117"; Value=3111},
    @{Sheet="CUL"; Cell="L117"; Value=9333},
    @{Sheet="CUL"; Cell="N117"; Value=-16217},
    @{Sheet="CUL"; Cell="H130"; Value=0},
    @{Sheet="CUL"; Cell="I130"; Value=0},
    @{Sheet="CUL"; Cell="J130"; Value=0},
    @{Sheet="CUL"; Cell="K130"; Value=0},
    @{Sheet="CUL"; Cell="L130"; Value=0},
    @{Sheet="CUL"; Cell="M130"; Value=$null},
    @{Sheet="CUL"; Cell="N130"; Value=$null},
    @{Sheet="CUL"; Cell="H132"; Value=1769},
    @{Sheet="CUL"; Cell="I132"; Value=1653.5},
    @{Sheet="CUL"; Cell="J132"; Value=2000},
    @{Sheet="CUL"; Cell="K132"; Value=14881.5},
    @{Sheet="CUL"; Cell="L132"; Value=18000},
    @{Sheet="CUL"; Cell="M132"; Value=-12351.5},
    @{Sheet="CUL"; Cell="N132"; Value=-23060},
    @{Sheet="GSM"; Cell="H42"; Value=65000},
    @{Sheet="GSM"; Cell="J42"; Value=65000},
    @{Sheet="GSM"; Cell="L42"; Value=65000},
    @{Sheet="GSM"; Cell="N42"; Value=-65970},
    @{Sheet="GSM"; Cell="H70"; Value=7499.4},
    @{Sheet="GSM"; Cell="J70"; Value=7499.4},
    @{Sheet="GSM"; Cell="L70"; Value=7499.4},
    @{Sheet="GSM"; Cell="N70"; Value=-8039.4},
    @{Sheet="GSM"; Cell="H73"; Value=7499.4},
    @{Sheet="GSM"; Cell="J73"; Value=7499.4},
    @{Sheet="GSM"; Cell="L73"; Value=7499.4},
    @{Sheet="GSM"; Cell="N73"; Value=-9371.4},
    @{Sheet="GSM"; Cell="H80"; Value=3924.75},
    @{Sheet="GSM"; Cell="I80"; Value=3499.5},
    @{Sheet="GSM"; Cell="J80"; Value=4350},
    @{Sheet="GSM"; Cell="K80"; Value=3499.5},
    @{Sheet="GSM"; Cell="L80"; Value=4350},
    @{Sheet="GSM"; Cell="M80"; Value=-2501.5},
    @{Sheet="GSM"; Cell="N80"; Value=-6346},
    @{Sheet="GSM"; Cell="H83"; Value=3924.75},
    @{Sheet="GSM"; Cell="I83"; Value=3499.5},
    @{Sheet="GSM"; Cell="J83"; Value=4350},
    @{Sheet="GSM"; Cell="K83"; Value=17497.5},
    @{Sheet="GSM"; Cell="L83"; Value=21750},
    @{Sheet="GSM"; Cell="M83"; Value=-12505.5},
    @{Sheet="GSM"; Cell="N83"; Value=-31734},
    @{Sheet="GSM"; Cell="H102"; Value=2472},
    @{Sheet="GSM"; Cell="J102"; Value=3006.2727},
    @{Sheet="GSM"; Cell="L102"; Value=3006.2727},
    @{Sheet="GSM"; Cell="N102"; Value=-6250.2727},
    @{Sheet="GSM"; Cell="H113"; Value=4023.3333},
    @{Sheet="GSM"; Cell="I113"; Value=2802.5},
    @{Sheet="GSM"; Cell="J113"; Value=5000},
    @{Sheet="GSM"; Cell="K113"; Value=2802.5},
    @{Sheet="GSM"; Cell="L113"; Value=5000},
    @{Sheet="GSM"; Cell="M113"; Value=-632.5},
    @{Sheet="GSM"; Cell="N113"; Value=-9340},
    @{Sheet="GSM"; Cell="H115"; Value=65000},
    @{Sheet="GSM"; Cell="J115"; Value=65000},
    @{Sheet="GSM"; Cell="L115"; Value=65000},
    @{Sheet="GSM"; Cell="N115"; Value=-67350},
    @{Sheet="GSM"; Cell="H126"; Value=4217.875},
    @{Sheet="GSM"; Cell="I126"; Value=3625},
    @{Sheet="GSM"; Cell="J126"; Value=4679},
    @{Sheet="GSM"; Cell="K126"; Value=10875},
    @{Sheet="GSM"; Cell="L126"; Value=14037},
    @{Sheet="GSM"; Cell="M126"; Value=-8405},
    @{Sheet="GSM"; Cell="N126"; Value=-18977},
    @{Sheet="GSM"; Cell="H132"; Value=2911.0557},
    @{Sheet="GSM"; Cell="I132"; Value=1376.4},
    @{Sheet="GSM"; Cell="K132"; Value=4129.200000000001},
    @{Sheet="GSM"; Cell="M132"; Value=-1599.200000000001},
    @{Sheet="LTW"; Cell="H16"; Value=990.1111},
    @{Sheet="LTW"; Cell="I16"; Value=990.1111},
    @{Sheet="LTW"; Cell="K16"; Value=990.1111},
    @{Sheet="LTW"; Cell="M16"; Value=-820.1111},
    @{Sheet="LTW"; Cell="H40"; Value=4780},
    @{Sheet="LTW"; Cell="I40"; Value=4450},
    @{Sheet="LTW"; Cell="K40"; Value=4450},
    @{Sheet="LTW"; Cell="M40"; Value=-4314},
    @{Sheet="LTW"; Cell="H46"; Value=3050.2942},
    @{Sheet="LTW"; Cell="J46"; Value=3345.3},
    @{Sheet="LTW"; Cell="L46"; Value=3345.3},
    @{Sheet="LTW"; Cell="N46"; Value=-3721.3},
    @{Sheet="LTW"; Cell="H100"; Value=1799.75},
    @{Sheet="LTW"; Cell="I100"; Value=1099.5},
    @{Sheet="LTW"; Cell="K100"; Value=1099.5},
    @{Sheet="LTW"; Cell="M100"; Value=-558.5},
    @{Sheet="LTW"; Cell="H132"; Value=5183.8184},
    @{Sheet="LTW"; Cell="I132"; Value=4767.8184},
    @{Sheet="LTW"; Cell="J132"; Value=5599.8184},
    @{Sheet="LTW"; Cell="K132"; Value=14303.4552},
    @{Sheet="LTW"; Cell="L132"; Value=16799.4552},
    @{Sheet="LTW"; Cell="M132"; Value=-11773.4552},
    @{Sheet="LTW"; Cell="N132"; Value=-21859.4552},
    @{Sheet="LTW"; Cell="H136"; Value=2171.2327},
    @{Sheet="LTW"; Cell="I136"; Value=2137.0571},
    @{Sheet="LTW"; Cell="J136"; Value=2320.75},
    @{Sheet="LTW"; Cell="K136"; Value=6411.1713},
    @{Sheet="LTW"; Cell="L136"; Value=6962.25},
    @{Sheet="LTW"; Cell="M136"; Value=-3861.1713},
    @{Sheet="LTW"; Cell="N136"; Value=-12062.25},
    @{Sheet="WVR"; Cell="H107"; Value=449.5},
    @{Sheet="WVR"; Cell="I107"; Value=449.5},
    @{Sheet="WVR"; Cell="K107"; Value=1348.5},
    @{Sheet="WVR"; Cell="M107"; Value=571.5},
    @{Sheet="WVR"; Cell="H132"; Value=620.2381},
    @{Sheet="WVR"; Cell="I132"; Value=623.9737},
    @{Sheet="WVR"; Cell="K132"; Value=1871.9211},
    @{Sheet="WVR"; Cell="M132"; Value=658.0789},
    @{Sheet="WVR"; Cell="H136"; Value=2824.22},
    @{Sheet="WVR"; Cell="I136"; Value=1433.6},
    @{Sheet="WVR"; Cell="K136"; Value=4300.799999999999},
    @{Sheet="WVR"; Cell="M136"; Value=-1750.799999999999}
)

$wb = $excel.ActiveWorkbook

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e.Sheet)
    $cell = $ws.Range($e.Cell)
    if ($e.Value -eq $null) {
        $cell.Value = ""
    } else {
        $cell.Value = $e.Value
    }
}

Write-Host "Applied $($edits.Count) cell updates across $(($edits | ForEach-Object { $_.Sheet } | Sort-Object -Unique).Count) sheets"
